$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.188.54"
$ws.Range("D2").Style = $style

$style = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("E2").Style = $style

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.931.03"
$ws.Range("D3").Style = $style

$style = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E3").Style = $style

$style = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E4").Style = $style

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.79"
$ws.Range("D5").Style = $style

$style = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E5").Style = $style

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.50"
$ws.Range("D6").Style = $style

$style = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E6").Style = $style

$style = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E7").Style = $style

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("D8").Style = $style

$style = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E8").Style = $style

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.04"
$ws.Range("D9").Style = $style

$style = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.37%  "
$ws.Range("E9").Style = $style

$style = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E10").Style = $style

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.441"
$ws.Range("D11").Style = $style

$style = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E11").Style = $style

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000226"
$ws.Range("D12").Style = $style

$style = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E12").Style = $style

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.78"
$ws.Range("D13").Style = $style

$style = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E13").Style = $style

$style = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E14").Style = $style

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.417.87"
$ws.Range("D15").Style = $style

$style = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("E15").Style = $style

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.216.18"
$ws.Range("D16").Style = $style

$style = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("E16").Style = $style

$style = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("E17").Style = $style

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.931.77"
$ws.Range("D18").Style = $style

$style = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E18").Style = $style

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "435.76"
$ws.Range("D19").Style = $style

$style = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("E19").Style = $style

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("D20").Style = $style

$style = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("E20").Style = $style

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.680"
$ws.Range("D21").Style = $style

$style = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E21").Style = $style

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.11"
$ws.Range("D22").Style = $style

$style = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E22").Style = $style

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.57"
$ws.Range("D23").Style = $style

$style = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("E23").Style = $style

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.06"
$ws.Range("D24").Style = $style

$style = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("E24").Style = $style

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("D25").Style = $style

$style = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E25").Style = $style

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("D26").Style = $style

$style = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E26").Style = $style

$style = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E27").Style = $style

$style = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E28").Style = $style

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.61"
$ws.Range("D29").Style = $style

$style = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E29").Style = $style

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("D30").Style = $style

$style = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E30").Style = $style

$style = $ws.Range("B31").Style
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("B31").Style = $style

$style = $ws.Range("C31").Style
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C31").Style = $style

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.80"
$ws.Range("D31").Style = $style

$style = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E31").Style = $style

$style = $ws.Range("B32").Style
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Hedera"
$ws.Range("B32").Style = $style

$style = $ws.Range("C32").Style
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C32").Style = $style

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.110"
$ws.Range("D32").Style = $style

$style = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("E32").Style = $style

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = $style

$style = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E33").Style = $style

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0869"
$ws.Range("D34").Style = $style

$style = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("E34").Style = $style

$style = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E35").Style = $style

$style = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E36").Style = $style

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.00"
$ws.Range("D37").Style = $style

$style = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E37").Style = $style

$style = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E38").Style = $style

$style = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E39").Style = $style

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.60"
$ws.Range("D40").Style = $style

$style = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E40").Style = $style

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.28"
$ws.Range("D41").Style = $style

$style = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.49%  "
$ws.Range("E41").Style = $style

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.285"
$ws.Range("D42").Style = $style

$style = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("E42").Style = $style

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "374.40"
$ws.Range("D43").Style = $style

$style = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E43").Style = $style

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0347"
$ws.Range("D44").Style = $style

$style = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E44").Style = $style

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.709.16"
$ws.Range("D45").Style = $style

$style = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("E45").Style = $style

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.95"
$ws.Range("D46").Style = $style

$style = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("E46").Style = $style

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.86"
$ws.Range("D48").Style = $style

$style = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E48").Style = $style

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.106"
$ws.Range("D49").Style = $style

$style = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("E49").Style = $style

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("D50").Style = $style

$style = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E50").Style = $style

$style = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.06%  "
$ws.Range("E51").Style = $style
